$d = $word.ActiveDocument

# Update the date heading paragraph (outside the table). Scope the Find to an
# explicit Document.Range(start, end) -- using a sub-objects .Range directly
# as the Find target does not reliably clip the search to that sub-range in
# this host, so we rebuild the range from $d.Range(...) before every Find.
$dateParagraph = $d.Paragraphs.Item(1)
$dateRange = $d.Range($dateParagraph.Range.Start, $dateParagraph.Range.End)
$dateRange.Find.Execute("2023-04-16 Sunday", $true, $false, $false, $false, $false, $true, 0, $false, "2023-04-17 Monday", 2) | Out-Null

# Update each arithmetic-problem cell in the table, in document order.
# Each cell is addressed via an explicit Document.Range(start, end) built from
# the cells own Range.Start/.End, and Find is scoped (Wrap = wdFindStop) so it
# cannot match text belonging to a different cell -- this matters because some
# expressions (e.g. "74-32=") occur in more than one cell with different
# replacements, and some old values are substrings of other cells old values
# (e.g. "7+20=" is contained in "77+20=").
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("11+74=", $true, $false, $false, $false, $false, $true, 0, $false, "45-22=", 2) | Out-Null
$cell = $t.Cell(1, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("64+18=", $true, $false, $false, $false, $false, $true, 0, $false, "25+3=", 2) | Out-Null
$cell = $t.Cell(1, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("62-28=", $true, $false, $false, $false, $false, $true, 0, $false, "40-14=", 2) | Out-Null
$cell = $t.Cell(1, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("2+26=", $true, $false, $false, $false, $false, $true, 0, $false, "85-16=", 2) | Out-Null
$cell = $t.Cell(1, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("17+55=", $true, $false, $false, $false, $false, $true, 0, $false, "55+5=", 2) | Out-Null
$cell = $t.Cell(2, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("4+69=", $true, $false, $false, $false, $false, $true, 0, $false, "38+51=", 2) | Out-Null
$cell = $t.Cell(2, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("31+49=", $true, $false, $false, $false, $false, $true, 0, $false, "5-1=", 2) | Out-Null
$cell = $t.Cell(2, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("5+86=", $true, $false, $false, $false, $false, $true, 0, $false, "26-22=", 2) | Out-Null
$cell = $t.Cell(2, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("2+52=", $true, $false, $false, $false, $false, $true, 0, $false, "7+2=", 2) | Out-Null
$cell = $t.Cell(2, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("82+0=", $true, $false, $false, $false, $false, $true, 0, $false, "17+82=", 2) | Out-Null
$cell = $t.Cell(3, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("10+16=", $true, $false, $false, $false, $false, $true, 0, $false, "29-7=", 2) | Out-Null
$cell = $t.Cell(3, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("18-6=", $true, $false, $false, $false, $false, $true, 0, $false, "8+45=", 2) | Out-Null
$cell = $t.Cell(3, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("49+21=", $true, $false, $false, $false, $false, $true, 0, $false, "51-15=", 2) | Out-Null
$cell = $t.Cell(3, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("5+55=", $true, $false, $false, $false, $false, $true, 0, $false, "50-41=", 2) | Out-Null
$cell = $t.Cell(3, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("51-32=", $true, $false, $false, $false, $false, $true, 0, $false, "63+17=", 2) | Out-Null
$cell = $t.Cell(4, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("84+9=", $true, $false, $false, $false, $false, $true, 0, $false, "31+64=", 2) | Out-Null
$cell = $t.Cell(4, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("16+74=", $true, $false, $false, $false, $false, $true, 0, $false, "90-40=", 2) | Out-Null
$cell = $t.Cell(4, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("93-0=", $true, $false, $false, $false, $false, $true, 0, $false, "88-53=", 2) | Out-Null
$cell = $t.Cell(4, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("75-27=", $true, $false, $false, $false, $false, $true, 0, $false, "51-18=", 2) | Out-Null
$cell = $t.Cell(4, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("73+22=", $true, $false, $false, $false, $false, $true, 0, $false, "34-4=", 2) | Out-Null
$cell = $t.Cell(5, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("0+2=", $true, $false, $false, $false, $false, $true, 0, $false, "8+34=", 2) | Out-Null
$cell = $t.Cell(5, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("19-1=", $true, $false, $false, $false, $false, $true, 0, $false, "65+13=", 2) | Out-Null
$cell = $t.Cell(5, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("17-11=", $true, $false, $false, $false, $false, $true, 0, $false, "48-36=", 2) | Out-Null
$cell = $t.Cell(5, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("85-7=", $true, $false, $false, $false, $false, $true, 0, $false, "26+54=", 2) | Out-Null
$cell = $t.Cell(5, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("48+31=", $true, $false, $false, $false, $false, $true, 0, $false, "26-13=", 2) | Out-Null
$cell = $t.Cell(6, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("14+45=", $true, $false, $false, $false, $false, $true, 0, $false, "39-1=", 2) | Out-Null
$cell = $t.Cell(6, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("7+20=", $true, $false, $false, $false, $false, $true, 0, $false, "34+29=", 2) | Out-Null
$cell = $t.Cell(6, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("62+23=", $true, $false, $false, $false, $false, $true, 0, $false, "12-0=", 2) | Out-Null
$cell = $t.Cell(6, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("14+24=", $true, $false, $false, $false, $false, $true, 0, $false, "68+2=", 2) | Out-Null
$cell = $t.Cell(6, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("53+15=", $true, $false, $false, $false, $false, $true, 0, $false, "78-62=", 2) | Out-Null
$cell = $t.Cell(7, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("46+44=", $true, $false, $false, $false, $false, $true, 0, $false, "68-45=", 2) | Out-Null
$cell = $t.Cell(7, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("26+28=", $true, $false, $false, $false, $false, $true, 0, $false, "91-70=", 2) | Out-Null
$cell = $t.Cell(7, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("19+49=", $true, $false, $false, $false, $false, $true, 0, $false, "26+3=", 2) | Out-Null
$cell = $t.Cell(7, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("35+5=", $true, $false, $false, $false, $false, $true, 0, $false, "97-32=", 2) | Out-Null
$cell = $t.Cell(7, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("59-18=", $true, $false, $false, $false, $false, $true, 0, $false, "67-59=", 2) | Out-Null
$cell = $t.Cell(8, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("3+13=", $true, $false, $false, $false, $false, $true, 0, $false, "36+43=", 2) | Out-Null
$cell = $t.Cell(8, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("83-40=", $true, $false, $false, $false, $false, $true, 0, $false, "24+69=", 2) | Out-Null
$cell = $t.Cell(8, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("93+1=", $true, $false, $false, $false, $false, $true, 0, $false, "84+1=", 2) | Out-Null
$cell = $t.Cell(8, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("87-75=", $true, $false, $false, $false, $false, $true, 0, $false, "92-38=", 2) | Out-Null
$cell = $t.Cell(8, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("58-7=", $true, $false, $false, $false, $false, $true, 0, $false, "35-12=", 2) | Out-Null
$cell = $t.Cell(9, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("32+32=", $true, $false, $false, $false, $false, $true, 0, $false, "30+44=", 2) | Out-Null
$cell = $t.Cell(9, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("82-54=", $true, $false, $false, $false, $false, $true, 0, $false, "34+40=", 2) | Out-Null
$cell = $t.Cell(9, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("81-80=", $true, $false, $false, $false, $false, $true, 0, $false, "96-78=", 2) | Out-Null
$cell = $t.Cell(9, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("74-32=", $true, $false, $false, $false, $false, $true, 0, $false, "75+6=", 2) | Out-Null
$cell = $t.Cell(9, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("16-12=", $true, $false, $false, $false, $false, $true, 0, $false, "72-31=", 2) | Out-Null
$cell = $t.Cell(10, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("4+95=", $true, $false, $false, $false, $false, $true, 0, $false, "78+17=", 2) | Out-Null
$cell = $t.Cell(10, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("54-9=", $true, $false, $false, $false, $false, $true, 0, $false, "4+89=", 2) | Out-Null
$cell = $t.Cell(10, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("67-49=", $true, $false, $false, $false, $false, $true, 0, $false, "38+14=", 2) | Out-Null
$cell = $t.Cell(10, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("96-0=", $true, $false, $false, $false, $false, $true, 0, $false, "59+7=", 2) | Out-Null
$cell = $t.Cell(10, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("73-23=", $true, $false, $false, $false, $false, $true, 0, $false, "89-15=", 2) | Out-Null
$cell = $t.Cell(11, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("28+11=", $true, $false, $false, $false, $false, $true, 0, $false, "51+28=", 2) | Out-Null
$cell = $t.Cell(11, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("86-7=", $true, $false, $false, $false, $false, $true, 0, $false, "65-47=", 2) | Out-Null
$cell = $t.Cell(11, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("48-43=", $true, $false, $false, $false, $false, $true, 0, $false, "28+30=", 2) | Out-Null
$cell = $t.Cell(11, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("32-3=", $true, $false, $false, $false, $false, $true, 0, $false, "4+6=", 2) | Out-Null
$cell = $t.Cell(11, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("24+13=", $true, $false, $false, $false, $false, $true, 0, $false, "47+52=", 2) | Out-Null
$cell = $t.Cell(12, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("30-5=", $true, $false, $false, $false, $false, $true, 0, $false, "36-14=", 2) | Out-Null
$cell = $t.Cell(12, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("84-47=", $true, $false, $false, $false, $false, $true, 0, $false, "14+43=", 2) | Out-Null
$cell = $t.Cell(12, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("74-32=", $true, $false, $false, $false, $false, $true, 0, $false, "3+0=", 2) | Out-Null
$cell = $t.Cell(12, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("64-20=", $true, $false, $false, $false, $false, $true, 0, $false, "33+4=", 2) | Out-Null
$cell = $t.Cell(12, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("82-56=", $true, $false, $false, $false, $false, $true, 0, $false, "10+12=", 2) | Out-Null
$cell = $t.Cell(13, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("61+27=", $true, $false, $false, $false, $false, $true, 0, $false, "80+15=", 2) | Out-Null
$cell = $t.Cell(13, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("61-52=", $true, $false, $false, $false, $false, $true, 0, $false, "47+46=", 2) | Out-Null
$cell = $t.Cell(13, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("50-27=", $true, $false, $false, $false, $false, $true, 0, $false, "9+71=", 2) | Out-Null
$cell = $t.Cell(13, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("39+46=", $true, $false, $false, $false, $false, $true, 0, $false, "63-17=", 2) | Out-Null
$cell = $t.Cell(13, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("90-50=", $true, $false, $false, $false, $false, $true, 0, $false, "66-46=", 2) | Out-Null
$cell = $t.Cell(14, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("79-72=", $true, $false, $false, $false, $false, $true, 0, $false, "38+35=", 2) | Out-Null
$cell = $t.Cell(14, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("59-52=", $true, $false, $false, $false, $false, $true, 0, $false, "53+23=", 2) | Out-Null
$cell = $t.Cell(14, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("19+48=", $true, $false, $false, $false, $false, $true, 0, $false, "5+19=", 2) | Out-Null
$cell = $t.Cell(14, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("17+63=", $true, $false, $false, $false, $false, $true, 0, $false, "52+31=", 2) | Out-Null
$cell = $t.Cell(14, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("22+3=", $true, $false, $false, $false, $false, $true, 0, $false, "94-35=", 2) | Out-Null
$cell = $t.Cell(15, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("49-32=", $true, $false, $false, $false, $false, $true, 0, $false, "9+1=", 2) | Out-Null
$cell = $t.Cell(15, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("35+57=", $true, $false, $false, $false, $false, $true, 0, $false, "31+51=", 2) | Out-Null
$cell = $t.Cell(15, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("67-66=", $true, $false, $false, $false, $false, $true, 0, $false, "35+53=", 2) | Out-Null
$cell = $t.Cell(15, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("23-5=", $true, $false, $false, $false, $false, $true, 0, $false, "15+80=", 2) | Out-Null
$cell = $t.Cell(15, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("77+20=", $true, $false, $false, $false, $false, $true, 0, $false, "7+4=", 2) | Out-Null
$cell = $t.Cell(16, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("52-42=", $true, $false, $false, $false, $false, $true, 0, $false, "1+63=", 2) | Out-Null
$cell = $t.Cell(16, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("70-70=", $true, $false, $false, $false, $false, $true, 0, $false, "69-53=", 2) | Out-Null
$cell = $t.Cell(16, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("58-31=", $true, $false, $false, $false, $false, $true, 0, $false, "96-76=", 2) | Out-Null
$cell = $t.Cell(16, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("91-89=", $true, $false, $false, $false, $false, $true, 0, $false, "54-26=", 2) | Out-Null
$cell = $t.Cell(16, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("2+37=", $true, $false, $false, $false, $false, $true, 0, $false, "19+77=", 2) | Out-Null
$cell = $t.Cell(17, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("54+39=", $true, $false, $false, $false, $false, $true, 0, $false, "84-35=", 2) | Out-Null
$cell = $t.Cell(17, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("92-33=", $true, $false, $false, $false, $false, $true, 0, $false, "26-22=", 2) | Out-Null
$cell = $t.Cell(17, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("68+7=", $true, $false, $false, $false, $false, $true, 0, $false, "86+11=", 2) | Out-Null
$cell = $t.Cell(17, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("93+4=", $true, $false, $false, $false, $false, $true, 0, $false, "66+4=", 2) | Out-Null
$cell = $t.Cell(17, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("72-22=", $true, $false, $false, $false, $false, $true, 0, $false, "31+47=", 2) | Out-Null
$cell = $t.Cell(18, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("82+17=", $true, $false, $false, $false, $false, $true, 0, $false, "44+42=", 2) | Out-Null
$cell = $t.Cell(18, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("27+47=", $true, $false, $false, $false, $false, $true, 0, $false, "48+39=", 2) | Out-Null
$cell = $t.Cell(18, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("72-50=", $true, $false, $false, $false, $false, $true, 0, $false, "2+68=", 2) | Out-Null
$cell = $t.Cell(18, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("91-61=", $true, $false, $false, $false, $false, $true, 0, $false, "70-16=", 2) | Out-Null
$cell = $t.Cell(18, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("36+5=", $true, $false, $false, $false, $false, $true, 0, $false, "57-22=", 2) | Out-Null
$cell = $t.Cell(19, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("15+48=", $true, $false, $false, $false, $false, $true, 0, $false, "23+11=", 2) | Out-Null
$cell = $t.Cell(19, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("93-3=", $true, $false, $false, $false, $false, $true, 0, $false, "95-17=", 2) | Out-Null
$cell = $t.Cell(19, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("27+53=", $true, $false, $false, $false, $false, $true, 0, $false, "69-41=", 2) | Out-Null
$cell = $t.Cell(19, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("47+29=", $true, $false, $false, $false, $false, $true, 0, $false, "81+12=", 2) | Out-Null
$cell = $t.Cell(19, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("97-75=", $true, $false, $false, $false, $false, $true, 0, $false, "99-85=", 2) | Out-Null
$cell = $t.Cell(20, 1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("81-3=", $true, $false, $false, $false, $false, $true, 0, $false, "43+12=", 2) | Out-Null
$cell = $t.Cell(20, 2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("52-6=", $true, $false, $false, $false, $false, $true, 0, $false, "19+36=", 2) | Out-Null
$cell = $t.Cell(20, 3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("51-42=", $true, $false, $false, $false, $false, $true, 0, $false, "12+7=", 2) | Out-Null
$cell = $t.Cell(20, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("56+32=", $true, $false, $false, $false, $false, $true, 0, $false, "63-21=", 2) | Out-Null
$cell = $t.Cell(20, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("88-50=", $true, $false, $false, $false, $false, $true, 0, $false, "26+44=", 2) | Out-Null
